# Generate Report for Handback
# This script updates the localization-status report to reflect a handback
# transform failure: the handback file name did not match the expected
# handoff-derived file name for both the zh-cn and de-de targets.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn = $wb.Worksheets.Item("zh-cn")
$dede = $wb.Worksheets.Item("de-de")

$newStatus = "Handback transform failed"

# Update the Status column (row 3 corresponds to the bb2ad15e... file) on
# every sheet that surfaces it.
$overview.Range("E3").Value = $newStatus
$overview.Range("F3").Value = $newStatus
$zhcn.Range("C3").Value = $newStatus
$dede.Range("C3").Value = $newStatus

# Widen the "Error Detail" column (P) on both locale sheets so the new
# diagnostic message is readable, and populate it with the mismatch detail.
$zhcn.Columns.Item(16).ColumnWidth = 39.14
$zhcn.Range("P3").Value = "Handback file name: gp1nuttt.hxa is different with handoff file name: bb2ad15e-9fbe-4506-8bd3-68063fccc924.98c64084ee5e4bc133a9d0b3f4a025d174e1a4d5.zh-cn."

$dede.Columns.Item(16).ColumnWidth = 39.14
$dede.Range("P3").Value = "Handback file name: gp1nuttt.hxa is different with handoff file name: bb2ad15e-9fbe-4506-8bd3-68063fccc924.98c64084ee5e4bc133a9d0b3f4a025d174e1a4d5.de-de."
